$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column for all existing rows (2-489)
# from 45202 (2023-10-03) to 45203 (2023-10-04).
$ws.Range("C2:C489").Value = 45203

# Row 489 gains an explicit row height (matches the new row that follows it).
$ws.Rows.Item(489).RowHeight = 15

# Append a new record as row 490.
$ws.Range("A490").Value = "A 46888-2023"

$ws.Range("B490").Value = 45201
$ws.Range("B490").NumberFormat = "YYYY-MM-DD"

$ws.Range("C490").Value = 45203
$ws.Range("C490").NumberFormat = "YYYY-MM-DD"

$ws.Range("D490").Value = "VÄSTERBOTTENS LÄN"
$ws.Range("E490").Value = "ROBERTSFORS"

$ws.Range("G490").Value = 1.7
$ws.Range("H490").Value = 0
$ws.Range("I490").Value = 0
$ws.Range("J490").Value = 0
$ws.Range("K490").Value = 0
$ws.Range("L490").Value = 0
$ws.Range("M490").Value = 0
$ws.Range("N490").Value = 0
$ws.Range("O490").Value = 0
$ws.Range("P490").Value = 0
$ws.Range("Q490").Value = 0

$ws.Range("R490").Value = ""
$ws.Range("R490").WrapText = $true
